# Kings Casino User Database - CRUD / reporting columns update
# Populates the "Unique ID" (A) and "Times Won" (I) columns for the
# existing user rows, and leaves the selection where the author last
# left off (J13) per the authoring session captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unique ID (column A) values for rows 2-11
$uniqueIds = @{
    2  = 48534
    3  = 3483423
    4  = 486343
    5  = 153485
    6  = 8646
    7  = 1237485
    8  = 486
    9  = 8643
    10 = 34856341
    11 = 748641
}

# Times Won (column I) values for rows 2-11
$timesWon = @{
    2  = 0
    3  = 5
    4  = 3
    5  = 4
    6  = 8
    7  = 3
    8  = 9
    9  = 4
    10 = 1
    11 = 3
}

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 1).Value = $uniqueIds[$row]
    $ws.Cells.Item($row, 9).Value = $timesWon[$row]
}

# Restore the author's last active selection
[void]$ws.Range("J13").Select()
